# Add two new slides ("Umsetzung" and "Ergebnis") after the existing
# "Aufgabe" slide (slide 3), matching its layout/design exactly by
# duplicating it and then swapping in the new title text.

$p = $ppt.ActivePresentation

# Slide 3 ("Aufgabe") is the template: same layout, same decorative
# rectangle, same body copy - only the title differs on the new slides.
$template = $p.Slides.Item(3)

# Duplicate() inserts the copy immediately after its source, so two
# duplicates in a row land at positions 4 and 5.
$template.Duplicate() | Out-Null
$template.Duplicate() | Out-Null

# Shape 2 on this layout is the title placeholder ("Titel 1").
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Text = "Umsetzung"
$p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "Ergebnis"

Write-Output "Slides after edit: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $title = $p.Slides.Item($i).Shapes.Item(2).TextFrame.TextRange.Text
    Write-Output "  $i : $title"
}
